$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.846.05"
$ws.Range("E2").Value = "  -1.75%  "

$ws.Range("D3").Value = "3.383.66"
$ws.Range("E3").Value = "  -0.91%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.85"
$ws.Range("E5").Value = "  -0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.62"
$ws.Range("E6").Value = "  -1.44%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "3.383.12"
$ws.Range("E8").Value = "  -0.91%  "

$ws.Range("E9").Value = "  -1.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.59"
$ws.Range("E10").Value = "  +1.16%  "

$ws.Range("E11").Value = "  -2.97%  "

$ws.Range("E12").Value = "  -1.40%  "

$ws.Range("D13").Value = "3.959.98"
$ws.Range("E13").Value = "  -1.04%  "

$ws.Range("E14").Value = "  +0.83%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.23"
$ws.Range("E15").Value = "  +3.02%  "

$ws.Range("E16").Value = "  -3.26%  "

$ws.Range("D17").Value = "3.382.22"
$ws.Range("E17").Value = "  -0.97%  "

$ws.Range("D18").Value = "60.955.76"
$ws.Range("E18").Value = "  -1.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.06"
$ws.Range("E19").Value = "  -0.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.50"
$ws.Range("E21").Value = "  -0.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "377.90"
$ws.Range("E22").Value = "  -3.17%  "

$ws.Range("E23").Value = "  -2.70%  "

$ws.Range("D24").Value = "3.525.09"
$ws.Range("E24").Value = "  -0.78%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.30"
$ws.Range("E26").Value = "  -0.29%  "

$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000125"
$ws.Range("E27").Value = "  -2.78%  "

$ws.Range("E28").Value = "  +11.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.171"
$ws.Range("E29").Value = "  +6.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.54"
$ws.Range("E30").Value = "  -2.14%  "

$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.17"
$ws.Range("E32").Value = "  -1.55%  "

$ws.Range("E33").Value = "  -0.72%  "

$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.76"
$ws.Range("E35").Value = "  +0.85%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.21"
$ws.Range("E36").Value = "  -4.65%  "

$ws.Range("E37").Value = "  -2.05%  "

$ws.Range("E38").Value = "  -1.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.60"
$ws.Range("E39").Value = "  +0.72%  "

$ws.Range("E40").Value = "  -4.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.07%  "

$ws.Range("E42").Value = "  -2.47%  "

$ws.Range("E43").Value = "  -3.31%  "

$ws.Range("E44").Value = "  -1.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.57"
$ws.Range("E45").Value = "  -0.54%  "

$ws.Range("E46").Value = "  -2.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.14"
$ws.Range("E47").Value = "  -3.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.47"
$ws.Range("E48").Value = "  +1.47%  "

$ws.Range("D49").Value = "2.450.03"
$ws.Range("E49").Value = "  +2.84%  "

$ws.Range("E50").Value = "  -2.31%  "

$ws.Range("E51").Value = "  +4.71%  "
